$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix capitalisation of the "age class" description (row 7, column B)
$ws.Range("B7").Value = "Age class"

# Reflect the selection/scroll state captured in the saved file:
# the whole data range got selected and the view scrolled down a bit.
$ws.Activate()
$ws.Range("A1:B38").Select()
$excel.ActiveWindow.ScrollRow = 9
